$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f1ee4580ee0>),`n                ('model',`n                 BaggingClassifier(estimator=SVC(C=1, kernel='linear',`n                                                 random_state=42),`n                                   n_estimators=50, random_state=42))])"
$ws.Range("B2").Value = 0.7147058823529411
$ws.Range("C2").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f1ef09867f0>, 'scaler': MinMaxScaler(), 'model__n_estimators': 50, 'model__estimator__kernel': 'linear', 'model__estimator__class_weight': None, 'model__estimator__C': 1}"
$ws.Range("H2").Value = 0.6201231950985656
$ws.Range("I2").Value = 0.07268029030058115
$ws.Range("J2").Value = 0.4631653149138443
$ws.Range("K2").Value = 0.1257141433292182

# Row 3 (only CV stats change)
$ws.Range("H3").Value = 0.7159173591320369
$ws.Range("I3").Value = 0.08700100323962788
$ws.Range("J3").Value = 0.6020948757125227
$ws.Range("K3").Value = 0.1233617220049524

# Row 4 (only CV stats change)
$ws.Range("H4").Value = 0.6158382192027458
$ws.Range("I4").Value = 0.1293650174961552
$ws.Range("J4").Value = 0.5207446347769877
$ws.Range("K4").Value = 0.1658447475079126

# Row 5
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f1ee4580c40>),`n                ('model',`n                 BaggingClassifier(estimator=SVC(C=1, kernel='linear',`n                                                 random_state=42),`n                                   n_estimators=5, random_state=42))])"
$ws.Range("B5").Value = 0.7823529411764707
$ws.Range("C5").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f1ef0975b80>, 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__kernel': 'linear', 'model__estimator__class_weight': None, 'model__estimator__C': 1}"
$ws.Range("H5").Value = 0.6925827286551395
$ws.Range("I5").Value = 0.1156594433779685
$ws.Range("J5").Value = 0.6064047638635874
$ws.Range("K5").Value = 0.1747704318057565
